$d = $word.ActiveDocument

# --- Paragraph 2: main templated revocation text ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.MoveEnd(1, -1)
$r2.Text = "I, {{ user.name.full(middle=’full’) }} of {{ user.address.line_one(bare=`"True`") }}, {{ user.address.city }}, {{ end_in_county(user.address.county) }}, {{ user.address.state }}, revoke the Durable Power of Attorney for Property dated {{ property_agent_date }}, empowering {{ property_agent.name.full(middle= ‘full’) }} to act as my agent. {% if any_property_successors == True %}In this document, the following successor agent(s) were named: {{ property_successors.full_names()}}. {% if property_replace_agent == True %} {{ new_property_agent.name.full(middle=’full’) }} shall now take the place of {{ property_agent.name.full(middle=’full’) }} as my agent for Power of Attorney for Property.{% endif %}{% for person in property_successors %}{% if person != property_who_is_promoted %}{% if property_replace_agent == True and person.remain == True %} {{ person.name.full(middle = ‘full’) }} shall remain a successor agent for Power of Attorney for Property.{% endif %}{% if property_replace_agent == True and person.remain == False %} I hereby revoke and withdraw all power and authority granted to {{ person.name.full(middle=’full’) }}.{% endif %}{% endif %}{% endfor %}{% if property_replace_agent == False %}I hereby revoke and withdraw all power and authority granted under the aforementioned Durable Power of Attorney for Property. {% endif %}{% endif %}"
$ins = $p2.Range
$ins.MoveEnd(1, -1)
$ins.Collapse(0)
$ins.InsertBreak(6)
$ins2 = $p2.Range
$ins2.MoveEnd(1, -1)
$ins2.Collapse(0)
$ins2.InsertAfter("Dated: _____________________")

# --- Paragraph 15: notary acknowledgement text ---
$p15 = $d.Paragraphs.Item(15)
$r15 = $p15.Range
$r15.MoveEnd(1, -1)
$r15.Text = "On this _______ day of _____________, _______, before me, ___________________________,  a notary public in said state, personally appeared ____________________________, personally known to me (or proved to me on the basis of satisfactory evidence) to be the person whose name is subscribed to the within instrument, and acknowledged to me that she/he executed the same in her/his authorized capacity, and that by her/his signature on the instrument, the person, or the entity upon behalf of which the person acted, executed the instrument. WITNESS my hand and official seal. "

# --- Move the _GoBack bookmark from the end of the document into the notary paragraph ---
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()
$bmPos = $p15.Range.Start + 75
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
